$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 794
$ws.Cells.Item(3, 6).Value = 2864
$ws.Cells.Item(6, 6).Value = 596
$ws.Cells.Item(11, 6).Value = 11893
$ws.Cells.Item(12, 6).Value = 6740
$ws.Cells.Item(30, 6).Value = 236
$ws.Cells.Item(32, 6).Value = 37
$ws.Cells.Item(34, 6).Value = 5056
$ws.Cells.Item(38, 6).Value = 693
$ws.Cells.Item(39, 6).Value = 1219

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 120
$ws.Cells.Item(11, 6).Value = 86

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 9102

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 9102
$ws.Cells.Item(5, 6).Value = 794
$ws.Cells.Item(6, 6).Value = 2864
$ws.Cells.Item(11, 6).Value = 596
$ws.Cells.Item(15, 6).Value = 11893
$ws.Cells.Item(16, 6).Value = 6740
$ws.Cells.Item(33, 6).Value = 236
$ws.Cells.Item(41, 6).Value = 1219
